# Informe-01-010013-A-TC-TM-TP.xlsx
#
# #8  Mejorar la generacion de SKOS Concept Schemes
# #16 Incluir descripciones para algunas medidas en los DSDs
# #17 Referenciada codelist que luego no tiene valores
# #19 Anadir propiedad en el DSD que identifique el ambito territorial aplicable
# #20 Generacion erronea de medidas en 01-080101-010105TC
#
# Row 1: column headers -> human readable Spanish labels (instead of the
#        raw dash-separated slugs).
# Row 2: DSD component identifiers, now including a territorial-scope
#        ("refArea") property per column, and an explicit measure
#        description (iaest-measure:numero-de-edificios) instead of the
#        iaest-dimension one.
# Row 3: component kind ("dim"/"medida") per column, following row 2.
# Row 4: datatype / concept-scheme / code-URI per column, including the
#        previously-missing "URI-Comunidad" territorial reference and
#        fixing the codelist column order.
# Row 5: the external mapping-file reference moves from column A to
#        column D (under "Clase de propietario").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers ---------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Número de edificios"
$ws.Cells.Item(1, 2).Value = "Comarca nombre"
$ws.Cells.Item(1, 3).Value = "Comarca código"
$ws.Cells.Item(1, 4).Value = "Clase de propietario"
$ws.Cells.Item(1, 5).Value = "Provincia código"
$ws.Cells.Item(1, 6).Value = "Aragón"
$ws.Cells.Item(1, 7).Value = "Municipio código"
$ws.Cells.Item(1, 8).Value = "Provincia nombre"
$ws.Cells.Item(1, 9).Value = "Municipio nombre"

# --- Row 2: DSD component identifiers ---------------------------------
$ws.Cells.Item(2, 1).Value = "iaest-measure:numero-de-edificios"
$ws.Cells.Item(2, 2).Value = "sdmx-dimension:refArea"
$ws.Cells.Item(2, 3).Value = "null"
$ws.Cells.Item(2, 4).Value = "iaest-dimension:clase-de-propietario"
$ws.Cells.Item(2, 5).Value = "null"
$ws.Cells.Item(2, 6).Value = "sdmx-dimension:refArea"
$ws.Cells.Item(2, 7).Value = "null"
$ws.Cells.Item(2, 8).Value = "sdmx-dimension:refArea"
$ws.Cells.Item(2, 9).Value = "sdmx-dimension:refArea"

# --- Row 3: component kind --------------------------------------------
$ws.Cells.Item(3, 1).Value = "medida"
$ws.Cells.Item(3, 2).Value = "dim"
$ws.Cells.Item(3, 3).Value = "null"
$ws.Cells.Item(3, 4).Value = "dim"
$ws.Cells.Item(3, 5).Value = "null"
$ws.Cells.Item(3, 6).Value = "dim"
$ws.Cells.Item(3, 7).Value = "null"
$ws.Cells.Item(3, 8).Value = "dim"
$ws.Cells.Item(3, 9).Value = "dim"

# --- Row 4: datatype / concept scheme / code URI ----------------------
$ws.Cells.Item(4, 1).Value = "xsd:int"
$ws.Cells.Item(4, 2).Value = "URI-comarca"
$ws.Cells.Item(4, 3).Value = "null"
$ws.Cells.Item(4, 4).Value = "skos:Concept"
$ws.Cells.Item(4, 5).Value = "null"
$ws.Cells.Item(4, 6).Value = "URI-Comunidad"
$ws.Cells.Item(4, 7).Value = "null"
$ws.Cells.Item(4, 8).Value = "URI-Provincia"
$ws.Cells.Item(4, 9).Value = "URI-Municipio"

# --- Row 5: mapping file reference moves from A5 to D5 -----------------
# Give D5 the same formatting A5 had (reuses the existing style record
# instead of minting a new one), then move the value over and wipe A5
# (content + formatting) so it goes back to a blank cell.
$ws.Cells.Item(5, 1).Copy()
$ws.Cells.Item(5, 4).PasteSpecial(-4122)
$ws.Cells.Item(5, 4).Value = "mapping-clase-de-propietario.xlsx"
$ws.Cells.Item(5, 1).Clear()
